# Update the question text on slide 4 ("Ribbon & Toolbar") to also ask
# about the Toolbar, per the commit diff:
#   "Was ist Ribbon überhaupt?"
#   -> "Was ist Ribbon überhaupt? Was hat es mit Toolbar zutun?"

$oldText = "Was ist Ribbon überhaupt?"
$newText = "Was ist Ribbon überhaupt? Was hat es mit Toolbar zutun?"

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(4)

function Find-ShapeWithText($shapes, $text) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq $text) {
            return $shp
        }
        if ($shp.Type -eq 6) {
            # msoGroup - the text box we need lives inside a group shape.
            $found = Find-ShapeWithText $shp.GroupItems $text
            if ($found -ne $null) { return $found }
        }
    }
    return $null
}

$target = Find-ShapeWithText $slide.Shapes $oldText
if ($target -ne $null) {
    $target.TextFrame.TextRange.Text = $newText
}
